$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update F column (initialValues) for rows 2-44: 1 -> 0
for ($r = 2; $r -le 44; $r++) {
    $ws.Cells.Item($r, 6).Value = 0
}

# Update E (randomWalkVariance) and F (initialValues) for rows 45-67
for ($r = 45; $r -le 67; $r++) {
    $ws.Cells.Item($r, 5).Value = 0.05
    $ws.Cells.Item($r, 6).Value = 0.5
}

# Update sheet view: topLeftCell and selection
$ws.Application.ActiveWindow.ScrollRow = 31
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("N36").Select()
